# Apply F-column (想去人数) updates across sheets per commit 456a3b4
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 894
$ws.Range("F3").Value = 13954
$ws.Range("F4").Value = 13782
$ws.Range("F5").Value = 1063
$ws.Range("F6").Value = 815
$ws.Range("F8").Value = 614
$ws.Range("F10").Value = 29
$ws.Range("F11").Value = 67
$ws.Range("F12").Value = 782
$ws.Range("F14").Value = 141
$ws.Range("F15").Value = 104
$ws.Range("F16").Value = 88
$ws.Range("F17").Value = 151
$ws.Range("F19").Value = 554
$ws.Range("F20").Value = 444
$ws.Range("F21").Value = 458
$ws.Range("F22").Value = 336
$ws.Range("F23").Value = 10
$ws.Range("F24").Value = 285
$ws.Range("F25").Value = 853
$ws.Range("F26").Value = 121
$ws.Range("F27").Value = 36
$ws.Range("F28").Value = 6
$ws.Range("F31").Value = 11
$ws.Range("F32").Value = 11

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 83
$ws.Range("F7").Value = 179
$ws.Range("F8").Value = 1608
$ws.Range("F13").Value = 79
$ws.Range("F15").Value = 1613

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 18
$ws.Range("F4").Value = 128

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 894
$ws.Range("F4").Value = 13954
$ws.Range("F5").Value = 13782
$ws.Range("F6").Value = 1063
$ws.Range("F7").Value = 815
$ws.Range("F9").Value = 614
$ws.Range("F11").Value = 29
$ws.Range("F12").Value = 67
$ws.Range("F13").Value = 782
$ws.Range("F17").Value = 18
$ws.Range("F18").Value = 141
$ws.Range("F19").Value = 104
$ws.Range("F20").Value = 88
$ws.Range("F21").Value = 151
$ws.Range("F24").Value = 83
$ws.Range("F25").Value = 128
$ws.Range("F26").Value = 554
$ws.Range("F27").Value = 444
$ws.Range("F28").Value = 458
$ws.Range("F29").Value = 336
$ws.Range("F30").Value = 10
$ws.Range("F31").Value = 285
$ws.Range("F32").Value = 853
$ws.Range("F33").Value = 179
$ws.Range("F34").Value = 1608
$ws.Range("F39").Value = 121
$ws.Range("F40").Value = 36
$ws.Range("F41").Value = 6
$ws.Range("F43").Value = 79
$ws.Range("F46").Value = 11
$ws.Range("F47").Value = 11
$ws.Range("F48").Value = 1613
